# Applies the cryptocurrency price/volume update described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "60.885.53"
$ws.Cells.Item(2, 5).Value = "  +2.54%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.612.37"
$ws.Cells.Item(3, 5).Value = "  +0.86%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "573.66"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.32%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "143.13"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.55%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.06%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.75%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.640.62"
$ws.Cells.Item(9, 5).Value = "  +1.53%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "6.54"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.19%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +2.76%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -1.70%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +7.16%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "3.078.96"
$ws.Cells.Item(14, 5).Value = "  +0.97%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "60.888.65"
$ws.Cells.Item(15, 5).Value = "  +2.58%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "23.64"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +4.79%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +2.52%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.626.38"
$ws.Cells.Item(18, 5).Value = "  +1.10%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +3.77%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.26"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +9.58%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "349.58"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +4.03%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "7.11"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +14.48%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.28%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.516"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +12.57%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.79%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.995"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.35%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "Kaspa"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.161"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.24%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +6.30%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "0.0₃0799"
$ws.Cells.Item(29, 5).Value = "  +2.00%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +11.15%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.11%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "6.30"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +3.20%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "161.71"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.88%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "19.54"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.51%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.29"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +5.35%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.962"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +8.93%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +4.11%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +6.87%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "37.70"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.40%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.860"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.84%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +3.42%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "297.92"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.76%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "139.07"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +10.56%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.995"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.31%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Stellar"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0987"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.63%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Hedera"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0553"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +2.48%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Mantle"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.607"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.76%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +3.65%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.60%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "19.68"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +5.69%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "4.82"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +6.50%  "
